# JS-Frameworks-Self-Evaluation-Protocol.xlsx
# Commit: "add issue, add project via modal added"
#
# Fills in scores / comments for the "Add project via modal" related
# rows (Basic Options section) and updates a couple of comment notes,
# then moves the viewport/selection like a reviewer scrolling further
# down the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Basic Options section: fill in the previously-empty "Score" (C) values
#     and the matching reviewer comments (E), row by row ---
$ws.Range("C13").Value = 5
$ws.Range("E13").Value = "done"

$ws.Range("C14").Value = 5
# E14 gets a new comment AND switches to left-aligned text (previously
# blank/centered numeric style) to better display the longer note.
$ws.Range("E14").Value = "evetualy todo: in case of error - to refresh so inputs to be empty"
$ws.Range("E14").HorizontalAlignment = -4131  # xlLeft

# Existing comment on row 28 (Add Project) gets replaced with a bug note.
$ws.Range("E28").Value = "Add project button not working"

$ws.Range("C15").Value = 10
$ws.Range("E15").Value = "links to project-view, issue-view "

$ws.Range("C16").Value = 5
$ws.Range("C17").Value = 5
$ws.Range("C21").Value = 10
$ws.Range("C22").Value = 5

# --- Advanced Options section ---
$ws.Range("C25").Value = 20
$ws.Range("E25").Value = "comments shown and add abiliti"

# --- Reviewer scrolled down and selected a different cell ---
$ws.Range("C26").Select()

$wb.Save()
